$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1) Rename run ids in rows 3 and 4
$ws.Range("A3").Value = "resnet50_002"
$ws.Range("A4").Value = "resnet50_003"

# 2) Fill in new row 6 with a new run entry (columns A-L only for now;
#    M/N are populated after the new table column is inserted below)
$ws.Range("A6").Value = "resnet50_004"
$ws.Range("B6").Value = "resnet50"
$ws.Range("C6").Value = 32
$ws.Range("D6").Value = "strict_full_balanced"
$ws.Range("E6").Value = "0,1,2,3,4,5,6,7"
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = "FrontViewDataset"
$ws.Range("L6").Value = "create_transform((3,320,224))"

# 3) Insert a new table column at the end (will become column N),
#    then move the "image dir" data (currently column M) into it:
#      - Copy(Destination) for the data rows so formatting/style travels too
#      - direct .Value assignment for the header so the table's column
#        name metadata is refreshed correctly
#    Finally repurpose column M as the new "num_epoch_unfreeze" column.
$newCol = $tbl.ListColumns.Add(14)

$ws.Range("M2:M30").Copy($ws.Range("N2:N30"))
$ws.Range("N1").Value = "image dir"

$ws.Range("M2:M30").ClearContents()
$ws.Range("M1").Value = "num_epoch_unfreeze"

# 4) Now that column N exists, finish populating row 6
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = "C:\Users\Daniel\Documents\Data\Batch1"

# 5) Update selection to match target view state
$ws.Range("M2").Select()
